$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: 290.000đ -> 270.000đ
$ws.Range("D16").Value = "270.000đ"
$ws.Range("E16").Value = "270.000đ"

# Row 21: 6.790.000đ -> 6.290.000đ
$ws.Range("D21").Value = "6.290.000đ"
$ws.Range("E21").Value = "6.290.000đ"

# Row 22: 10.990.000đ -> 10.490.000đ
$ws.Range("D22").Value = "10.490.000đ"
$ws.Range("E22").Value = "10.490.000đ"

# Row 26: 419.000đ -> 390.000đ
$ws.Range("D26").Value = "390.000đ"
$ws.Range("E26").Value = "390.000đ"
